$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New news items (title, link, abstract) for rows 2-11, reflecting the
# latest news feed. Row 4 (Inter Risk Services) stays unchanged.
$news = @(
    @{ Row = 2;  Title = "Alessandro Octaviani toma posse como superintendente da Susep";
       Link = "https://www.revistaapolice.com.br/2023/04/alessandro-octaviani-toma-posse-como-superintendente-da-susep/";
       Abstract = "O novo superintendente entrará em exercício na próxima segunda-feira, 10 de abril, dedicando a primeira semana a reuniões e despachos internos, para apresentação dos projetos em curso na autarquia" },

    @{ Row = 3;  Title = "Marcelo Blay e Manes Erlichman assumem nova posição na Creditas";
       Link = "https://www.revistaapolice.com.br/2023/04/marcelo-blay-e-manes-erlichman-assumem-nova-posicao-na-creditas/";
       Abstract = "Marcelo Blay e Manes Erlichman assumem a posição de senior advisors na Creditas, que adquiriu a Minuto Seguros em julho de 2021" },

    @{ Row = 5;  Title = "Crescimento de roubo e furto de veículos reforça importância do seguro automóvel";
       Link = "https://www.revistaapolice.com.br/2023/04/crescimento-de-roubo-e-furto-de-veiculos-reforca-importancia-do-seguro-automovel/";
       Abstract = "De acordo com dados do IBGE, o Brasil registrou cerca de 64 roubos e furtos de veículos por hora em 2022" },

    @{ Row = 6;  Title = "Susep participa de painéis na Oficina FIDES";
       Link = "https://www.revistaapolice.com.br/2023/04/susep-participa-de-paineis-na-oficina-fides/";
       Abstract = "Evento é realizado pela CNseg nos dias 13 e 14 de abril, em Brasília, abordando os temas Sustentabilidade e Riscos Cibernéticos" },

    @{ Row = 7;  Title = "Seguros SURA registra crescimento em Minas Gerais e no Centro-Oeste";
       Link = "https://www.revistaapolice.com.br/2023/04/seguros-sura-registra-crescimento-em-minas-gerais-e-no-centro-oeste/";
       Abstract = "Objetivo da seguradora é continuar desenvolvendo produtos para ampliar sua carteira de clientes nas regiões" },

    @{ Row = 8;  Title = "Seguros Unimed investe em data lakehouse com tecnologia da AWS";
       Link = "https://www.revistaapolice.com.br/2023/04/seguros-unimed-investe-em-data-lakehouse-com-tecnologia-da-aws/";
       Abstract = "Objetivo principal da seguradora é ser uma empresa data driven, garantindo o cumprimento das normas da LGPD" },

    @{ Row = 9;  Title = "Bradesco Vida e Previdência lança ‘Pensão Temporária Bradesco’";
       Link = "https://www.revistaapolice.com.br/2023/04/bradesco-vida-e-previdencia-lanca-pensao-temporaria-bradesco/";
       Abstract = "Novo produto oferece flexibilidade na escolha de beneficiários, pensão por até 20 anos e pecúlio para uso imediato em caso da falta do contratante" },

    @{ Row = 10; Title = "Open Insurance: Segunda fase gera grande expectativas para insurtechs";
       Link = "https://www.revistaapolice.com.br/2023/04/open-insurance-segunda-fase-gera-grande-expectativas-para-insurtechs/";
       Abstract = "Startups do setor de seguros têm mais facilidade para aderir ao novo sistema e explorar oportunidades com decisões baseadas em dados" },

    @{ Row = 11; Title = "Aumento de roubo de celulares reforça importância do seguro";
       Link = "https://www.revistaapolice.com.br/2023/04/aumento-de-roubo-de-celulares-reforca-importancia-do-seguro/";
       Abstract = "Somente em São Paulo, mais de duzentos mil celulares foram roubados em 2022 segundo dados da Secretária de Segurança Pública" }
)

foreach ($item in $news) {
    $r = $item.Row
    $ws.Range("B$r").Value = $item.Title
    $ws.Range("C$r").Value = $item.Link
    $ws.Range("D$r").Value = $item.Abstract
}
